$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for row 117 ("Apio", Primera,
# fecha 2021-09-16 => serial 44455). Insert a fresh row at 117 (pushing the
# existing rows 117-173 down to 118-174, which matches how every row's data
# shifted by one position in the diff) and populate it with the new reading.
$ws.Rows(117).Insert()

$ws.Range("A117").Value = 10
$ws.Range("B117").Value = "Vega Modelo de Temuco"
$ws.Range("C117").Value = "La Araucanía"
$ws.Range("D117").Value = 44455
$ws.Range("E117").Value = 9
$ws.Range("F117").Value = 100112017
$ws.Range("G117").Value = "Apio"
$ws.Range("H117").Value = "Americana (o)"
$ws.Range("I117").Value = "Primera"
$ws.Range("J117").Value = 80
$ws.Range("K117").Value = 9000
$ws.Range("L117").Value = 9000
$ws.Range("M117").Value = 9000
$ws.Range("N117").Value = "`$/docena de matas"
$ws.Range("O117").Value = "Provincia del Elquí"
$ws.Range("P117").Value = 1500
$ws.Range("Q117").Value = 6
$ws.Range("R117").Value = "Hortaliza"
